$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.538.06"
$ws.Range("E2").Value = "  -2.13%  "

$ws.Range("D3").Value = "2.344.51"
$ws.Range("E3").Value = "  -2.79%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "'319.88"
$ws.Range("E5").Value = "  -1.48%  "

$ws.Range("D6").Value = "'105.04"
$ws.Range("E6").Value = "  +1.07%  "

$ws.Range("E7").Value = "  -1.47%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.615"
$ws.Range("E9").Value = "  -6.46%  "

$ws.Range("D10").Value = "'40.72"
$ws.Range("E10").Value = "  -2.40%  "

$ws.Range("D11").Value = "'0.0922"
$ws.Range("E11").Value = "  -2.93%  "

$ws.Range("D12").Value = "'8.36"
$ws.Range("E12").Value = "  -2.99%  "

$ws.Range("D13").Value = "'0.988"
$ws.Range("E13").Value = "  -4.76%  "

$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").Value = "'15.93"
$ws.Range("E15").Value = "  -7.94%  "

$ws.Range("D16").Value = "2.698.79"
$ws.Range("E16").Value = "  -2.90%  "

$ws.Range("D17").Value = "2.390.61"
$ws.Range("E17").Value = "  -3.85%  "

$ws.Range("D18").Value = "42.660.72"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("D19").Value = "'7.71"
$ws.Range("E19").Value = "  +3.40%  "

$ws.Range("E20").Value = "  -4.08%  "

$ws.Range("D21").Value = "'77.67"
$ws.Range("E21").Value = "  +2.95%  "

$ws.Range("D22").Value = "'3.56"
$ws.Range("E22").Value = "  +1.51%  "

$ws.Range("D23").Value = "'260.72"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("E24").Value = "  -4.99%  "

$ws.Range("D25").Value = "'9.62"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("D27").Value = "'11.38"
$ws.Range("E27").Value = "  -5.08%  "

$ws.Range("D28").Value = "'23.21"
$ws.Range("E28").Value = "  +1.50%  "

$ws.Range("E29").Value = "  -1.15%  "

$ws.Range("D30").Value = "'175.14"
$ws.Range("E30").Value = "  -2.46%  "

$ws.Range("D31").Value = "'36.36"
$ws.Range("E31").Value = "  -4.28%  "

$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "'3.01"
$ws.Range("E32").Value = "  -6.86%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0890"
$ws.Range("E33").Value = "  -4.54%  "

$ws.Range("E34").Value = "  +2.09%  "

$ws.Range("E35").Value = "  -1.91%  "

$ws.Range("E36").Value = "  +5.47%  "

$ws.Range("D37").Value = "'4.59"
$ws.Range("E37").Value = "  -5.60%  "

$ws.Range("E38").Value = "  -4.24%  "

$ws.Range("D39").Value = "'3.76"
$ws.Range("E39").Value = "  -5.20%  "

$ws.Range("D40").Value = "'2.65"
$ws.Range("E40").Value = "  -8.88%  "

$ws.Range("E41").Value = "  -10.58%  "

$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "'70.10"
$ws.Range("E42").Value = "  +1.11%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.231"
$ws.Range("E43").Value = "  -1.41%  "

$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").Value = "'114.66"
$ws.Range("E45").Value = "  -7.94%  "

$ws.Range("D46").Value = "'11.86"
$ws.Range("E46").Value = "  -5.82%  "

$ws.Range("E47").Value = "  -3.30%  "

$ws.Range("D48").Value = "'9.21"
$ws.Range("E48").Value = "  -3.55%  "

$ws.Range("D49").Value = "'84.58"
$ws.Range("E49").Value = "  +6.51%  "

$ws.Range("D50").Value = "'72.86"
$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("D51").Value = "'0.0998"
$ws.Range("E51").Value = "  -1.23%  "
